# Apply updated crypto market data (price + volume) per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.113.85'
$ws.Range('E2').Value = '  -2.78%  '
$ws.Range('D3').Value = '1.845.16'
$ws.Range('E3').Value = '  -2.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7003'
$ws.Range('E5').Value = '  -5.83%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '237.41'
$ws.Range('E6').Value = '  -2.30%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3030'
$ws.Range('E8').Value = '  -4.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07408'
$ws.Range('E9').Value = '  +2.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.34'
$ws.Range('E10').Value = '  -6.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08110'
$ws.Range('E11').Value = '  -2.87%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.848.26'
$ws.Range('E12').Value = '  -13.49%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7248'
$ws.Range('E13').Value = '  -4.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.207'
$ws.Range('E14').Value = '  -3.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.07'
$ws.Range('E15').Value = '  -3.69%  '
$ws.Range('D16').Value = '29.062.66'
$ws.Range('E16').Value = '  -3.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.784'
$ws.Range('E17').Value = '  -5.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '241.70'
$ws.Range('E18').Value = '  -2.94%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007660'
$ws.Range('E19').Value = '  -2.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.01'
$ws.Range('E20').Value = '  -4.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9995'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').Value = '2.078.47'
$ws.Range('E22').Value = '  -3.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9996'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.573'
$ws.Range('E24').Value = '  -5.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1471'
$ws.Range('E25').Value = '  -5.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.97'
$ws.Range('E26').Value = '  -2.47%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.938'
$ws.Range('E27').Value = '  -4.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.03'
$ws.Range('E28').Value = '  -3.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.934'
$ws.Range('E29').Value = '  -5.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.375'
$ws.Range('E30').Value = '  -8.63%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.450'
$ws.Range('E31').Value = '  -3.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.490'
$ws.Range('E32').Value = '  -3.03%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.005'
$ws.Range('E33').Value = '  -5.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05203'
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.183'
$ws.Range('E35').Value = '  -5.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7108'
$ws.Range('E36').Value = '  -6.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.001'
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.648'
$ws.Range('E38').Value = '  -2.17%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01869'
$ws.Range('E39').Value = '  -5.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.667'
$ws.Range('E40').Value = '  -3.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9065'
$ws.Range('E41').Value = '  +5.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4289'
$ws.Range('E42').Value = '  -5.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.901'
$ws.Range('E43').Value = '  -4.18%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '70.03'
$ws.Range('E44').Value = '  -4.03%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '1.048.46'
$ws.Range('E45').Value = '  -5.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9997'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.47'
$ws.Range('E47').Value = '  -3.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.755'
$ws.Range('E48').Value = '  -6.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.110'
$ws.Range('E49').Value = '  -6.82%  '
$ws.Range('D50').Value = '1.985.23'
$ws.Range('E50').Value = '  -4.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.188'
$ws.Range('E51').Value = '  -3.97%  '
